$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Dueño"
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = "Consulta clientes"

$ws.Range("C7").Value = "Se modifica de crear a crud"
$ws.Range("A7").Value = "Bodega"
$ws.Range("B7").Value = 4

$ws.Range("A8").Value = "Administrador"
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = "CRUD usuario"

$ws.Range("A9").Value = "Administrador"
$ws.Range("B9").Value = 24
$ws.Range("C9").Value = "Habilita cuenta"

$ws.Range("A10").Value = "Administrador"
$ws.Range("B10").Value = 25
$ws.Range("C10").Value = "Clasifica usuario"

$ws.Range("A11").Value = "Administrador"
$ws.Range("B11").Value = 26
$ws.Range("C11").Value = "Clasifica producto"

$ws.Columns.Item(1).ColumnWidth = 12.877604166666666

$ws.Range("D6").Select()
